# Update the training schedule data on Sheet1:
#   E6: 7 -> 6
#   G6: -2 -> -3
#   H6: 14 -> 13
# and leave the selection on E6 (the last cell the author edited), matching
# the saved sheetView state in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E6").Value = 6
$ws.Range("G6").Value = -3
$ws.Range("H6").Value = 13

$ws.Range("E6").Select()
